$wb = $excel.ActiveWorkbook

# Duplicate the "Spain" sheet (last sheet) to create the new "Turkey" sheet right after it
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $spain)

# The newly copied sheet becomes the active sheet, placed immediately after Spain
$turkey = $wb.ActiveSheet
$turkey.Name = "Turkey"

# Update the market name / NGC code for the Turkey test data
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3310/T3325/T3291"

# Column widths differ slightly from the copied "Spain" sheet (bestFit content changed)
$turkey.Columns.Item(2).ColumnWidth = 29.25
$turkey.Columns.Item(4).ColumnWidth = 21.42

# Rows 3 and 5 on "Spain" had a custom (taller) row height; Turkey uses the default height
$turkey.Rows.Item(3).AutoFit()
$turkey.Rows.Item(5).AutoFit()

# Insert an extra "PROFILE Communicator" row (between RS800 and Zettler Printer)
$turkey.Rows.Item(14).Insert()
$turkey.Range("A13").Copy($turkey.Range("A14"))
$turkey.Range("A14").Value = "PROFILE Communicator"

# Restore the selection on the new sheet
[void]$turkey.Range("A11").Select()

# The previously-active "Spain" sheet is no longer the selected tab;
# its selection becomes the full used range
[void]$spain.Range("A1:D16").Select()
[void]$turkey.Activate()
